$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the value from B16 ("      nodhiambo01") entirely, matching the
# diff which drops the <c r="B16"> cell (and its shared string) altogether.
$ws.Range("B16").ClearContents()
